$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for new columns
$ws.Range("K1").Value = "Ventral"
$ws.Range("L1").Value = "Dorsal"

# Formulas: average of left/right VS -> Ventral, average of left/right DS -> Dorsal
$lastRow = 141
$ws.Range("K2:K$lastRow").Formula = "=AVERAGE(G2,H2)"
$ws.Range("L2:L$lastRow").Formula = "=AVERAGE(I2,J2)"
